$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write `text` into `cellAddr` as a genuine text (shared-string) value
# without Excel's numeric auto-detection kicking in (which would store it as
# a number) and without leaving a "Number stored as text" quote-prefix style
# behind (which a leading apostrophe would create). We do this by building
# the text in a scratch cell via a text formula (so its result is unambiguous
# text), copying it, and pasting *values only* into the destination - this
# carries over the shared-string text type but none of the source/number
# formatting.
function Set-TextValue($cellAddr, $text) {
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.ClearContents()
}

Set-TextValue "A13" "118451"
Set-TextValue "B13" "1008617738"
Set-TextValue "C13" "17705210"
Set-TextValue "D13" "6001"

Set-TextValue "A14" "118451"
Set-TextValue "B14" "1008617741"
Set-TextValue "C14" "17705210"
Set-TextValue "D14" "6001"

Set-TextValue "A15" "118451"
Set-TextValue "B15" "1008617742"
Set-TextValue "C15" "17705211"
Set-TextValue "D15" "6004"

Set-TextValue "A16" "118452"
Set-TextValue "B16" "1008617743"
Set-TextValue "C16" "17705212"
Set-TextValue "D16" "6005"

Set-TextValue "A17" "118452"
Set-TextValue "B17" "1008617765"
Set-TextValue "C17" "17705235"
Set-TextValue "D17" "6004"
